# ajout de fichier de conversion Excel en xml
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the B1 ("LastName") / C1 ("FirstName") header values, and clear
# their bold/black-font style (s="1") so both fall back to the default
# cell style (s="0").
$ws.Range("B1:C1").ClearFormats()
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# Every student CNE value in column A (rows 2-11) becomes 18000031.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = 18000031
}

# Update the sheet view: drop the topLeftCell="B1" scroll offset (back to
# A1) and move the active selection from I10 to F12.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F12").Select()
